$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the "Espinaca" series. It belongs
# right before the current row 518 (it sits between the two 2021-08-17
# entries and the new 2023-01-13 entry, chronologically), so insert a new
# row there which pushes the existing rows 518-543 down to 519-544.
$ws.Rows.Item(518).Insert()

# Populate the newly inserted row 518 with the new record. All attributes
# mirror the row that used to be at 518 (same market/product/grade/prices),
# only the date (column D) differs - a new weekly observation.
$ws.Range("A518").Value = 9
$ws.Range("B518").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C518").Value = "Metropolitana"
$ws.Range("D518").Value = 44939
$ws.Range("E518").Value = 13
$ws.Range("F518").Value = 100112012
$ws.Range("G518").Value = "Espinaca"
$ws.Range("H518").Value = "Sin especificar"
$ws.Range("I518").Value = "Primera"
$ws.Range("J518").Value = 160
$ws.Range("K518").Value = 7000
$ws.Range("L518").Value = 8000
$ws.Range("M518").Value = 7500
$ws.Range("N518").Value = "$/cuna 10 kilos"
$ws.Range("O518").Value = "Provincia de Chacabuco"
$ws.Range("P518").Value = 750
$ws.Range("Q518").Value = 10
$ws.Range("R518").Value = "Hortaliza"
